# Apply "new data" update to VideoInfo sheet:
# Appends 45 new video rows (326-370) with Content id, Video title,
# Classification1, Classification2, and VideoDuration columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 326
$ws.Cells.Item(326, 1).Value = 'AxgzzJposVo'
$ws.Cells.Item(326, 2).Value = 'Elementary Row Operations, Row Echelon Form, and Reduced Row Echelon Form'
$ws.Cells.Item(326, 3).Value = 'Math'
$ws.Cells.Item(326, 4).Value = 'Education'
$ws.Cells.Item(326, 5).Value = '53:48'

# Row 327
$ws.Cells.Item(327, 1).Value = 'eIv8muz9Hsk'
$ws.Cells.Item(327, 2).Value = 'Row/Column Space, Linear Independence, and Rank of a Matrix'
$ws.Cells.Item(327, 3).Value = 'Math'
$ws.Cells.Item(327, 4).Value = 'Education'
$ws.Cells.Item(327, 5).Value = '38:58'

# Row 328
$ws.Cells.Item(328, 1).Value = 'haJVEtLN6-k'
$ws.Cells.Item(328, 2).Value = 'Scalar  Functions, Vector Functions, and Vector Derivatives'
$ws.Cells.Item(328, 3).Value = 'Math'
$ws.Cells.Item(328, 4).Value = 'Education'
$ws.Cells.Item(328, 5).Value = '22:54'

# Row 329
$ws.Cells.Item(329, 1).Value = 'MPcfaNIREN0'
$ws.Cells.Item(329, 2).Value = 'Parameterizing Curves'
$ws.Cells.Item(329, 3).Value = 'Math'
$ws.Cells.Item(329, 4).Value = 'Education'
$ws.Cells.Item(329, 5).Value = '14:58'

# Row 330
$ws.Cells.Item(330, 1).Value = 'HH367um_Aho'
$ws.Cells.Item(330, 2).Value = 'Tangent to a Curve'
$ws.Cells.Item(330, 3).Value = 'Math'
$ws.Cells.Item(330, 4).Value = 'Education'
$ws.Cells.Item(330, 5).Value = '13:57'

# Row 331
$ws.Cells.Item(331, 1).Value = 'FoiuvPkFppg'
$ws.Cells.Item(331, 2).Value = 'Arc Length (AKA Length of a Curve)'
$ws.Cells.Item(331, 3).Value = 'Math'
$ws.Cells.Item(331, 4).Value = 'Education'
$ws.Cells.Item(331, 5).Value = '25:40'

# Row 332
$ws.Cells.Item(332, 1).Value = 'obeu4B8mXuw'
$ws.Cells.Item(332, 2).Value = 'Gradient of a Function and the Directional Derivative'
$ws.Cells.Item(332, 3).Value = 'Math'
$ws.Cells.Item(332, 4).Value = 'Education'
$ws.Cells.Item(332, 5).Value = '26:57'

# Row 333
$ws.Cells.Item(333, 1).Value = '7GXbPYzW5JA'
$ws.Cells.Item(333, 2).Value = 'Fourier Series'
$ws.Cells.Item(333, 3).Value = 'Math'
$ws.Cells.Item(333, 4).Value = 'Education'
$ws.Cells.Item(333, 5).Value = '48:13'

# Row 334
$ws.Cells.Item(334, 1).Value = 'W30U_rcThLg'
$ws.Cells.Item(334, 2).Value = 'Discrete Fourier Transform'
$ws.Cells.Item(334, 3).Value = 'Math'
$ws.Cells.Item(334, 4).Value = 'Education'
$ws.Cells.Item(334, 5).Value = '1:22:07'

# Row 335
$ws.Cells.Item(335, 1).Value = 'yfsSDynscEs'
$ws.Cells.Item(335, 2).Value = 'Fast Fourier Transform'
$ws.Cells.Item(335, 3).Value = 'Math'
$ws.Cells.Item(335, 4).Value = 'Education'
$ws.Cells.Item(335, 5).Value = '48:22'

# Row 336
$ws.Cells.Item(336, 1).Value = 'bKwfnulkt2U'
$ws.Cells.Item(336, 2).Value = 'Replace Microsoft Surface Pen Battery'
$ws.Cells.Item(336, 3).Value = 'DIY'
$ws.Cells.Item(336, 4).Value = 'Other'
$ws.Cells.Item(336, 5).Value = '2:40'

# Row 337
$ws.Cells.Item(337, 1).Value = '5ROS2-4ShmI'
$ws.Cells.Item(337, 2).Value = 'Obtaining and Using the MatlabLum Repository'
$ws.Cells.Item(337, 3).Value = 'Matlab'
$ws.Cells.Item(337, 4).Value = 'Education'
$ws.Cells.Item(337, 5).Value = '13:40'

# Row 338
$ws.Cells.Item(338, 1).Value = 'NS2FI6vR3BY'
$ws.Cells.Item(338, 2).Value = 'The MNIST Database'
$ws.Cells.Item(338, 3).Value = 'AIML'
$ws.Cells.Item(338, 4).Value = 'Education'
$ws.Cells.Item(338, 5).Value = '20:27'

# Row 339
$ws.Cells.Item(339, 1).Value = 'eF0Zv-GPzH0'
$ws.Cells.Item(339, 2).Value = 'Planet Bike Superflash 65R Bike Tail Light'
$ws.Cells.Item(339, 3).Value = 'Review'
$ws.Cells.Item(339, 4).Value = 'Other'
$ws.Cells.Item(339, 5).Value = '7:26'

# Row 340
$ws.Cells.Item(340, 1).Value = 'k-mli8-04RQ'
$ws.Cells.Item(340, 2).Value = 'Time to Double for a First and Second Order System'
$ws.Cells.Item(340, 3).Value = 'Controls'
$ws.Cells.Item(340, 4).Value = 'Education'
$ws.Cells.Item(340, 5).Value = '40:41'

# Row 341
$ws.Cells.Item(341, 1).Value = 'b1uIEnjqcZM'
$ws.Cells.Item(341, 2).Value = 'How to Replace a Dryer Thermal Fuse'
$ws.Cells.Item(341, 3).Value = 'DIY'
$ws.Cells.Item(341, 4).Value = 'Other'
$ws.Cells.Item(341, 5).Value = '9:08'

# Row 342
$ws.Cells.Item(342, 1).Value = 'UkM1h7-URo8'
$ws.Cells.Item(342, 2).Value = 'Demolish Concrete Using Only a Sledgehammer'
$ws.Cells.Item(342, 3).Value = 'DIY'
$ws.Cells.Item(342, 4).Value = 'Other'
$ws.Cells.Item(342, 5).Value = '4:23'

# Row 343
$ws.Cells.Item(343, 1).Value = 'AOR2u3dwUNM'
$ws.Cells.Item(343, 2).Value = 'How to Wash and Dry a Down Jacket'
$ws.Cells.Item(343, 3).Value = 'DIY'
$ws.Cells.Item(343, 4).Value = 'Other'
$ws.Cells.Item(343, 5).Value = '5:49'

# Row 344
$ws.Cells.Item(344, 1).Value = 'nZ6LcTjtCbs'
$ws.Cells.Item(344, 2).Value = 'Monitor Wonâ€™t Turn On â€“ Solved!'
$ws.Cells.Item(344, 3).Value = 'DIY'
$ws.Cells.Item(344, 4).Value = 'Other'
$ws.Cells.Item(344, 5).Value = '11:55'

# Row 345
$ws.Cells.Item(345, 1).Value = 'XO8KvIoCNbE'
$ws.Cells.Item(345, 2).Value = 'Dimensionless Aerodynamic Coefficients'
$ws.Cells.Item(345, 3).Value = 'FlightMechanics'
$ws.Cells.Item(345, 4).Value = 'Education'
$ws.Cells.Item(345, 5).Value = '19:23'

# Row 346
$ws.Cells.Item(346, 1).Value = 'c3bRUDvLTS4'
$ws.Cells.Item(346, 2).Value = 'How to Remove the Cap on a Swiffer Wet Jet Bottle'
$ws.Cells.Item(346, 3).Value = 'DIY'
$ws.Cells.Item(346, 4).Value = 'Other'
$ws.Cells.Item(346, 5).Value = '4:06'

# Row 347
$ws.Cells.Item(347, 1).Value = 'QexBVGVM690'
$ws.Cells.Item(347, 2).Value = 'The Jacobian Matrix'
$ws.Cells.Item(347, 3).Value = 'Math'
$ws.Cells.Item(347, 4).Value = 'Education'
$ws.Cells.Item(347, 5).Value = '40:21'

# Row 348
$ws.Cells.Item(348, 1).Value = 'TULEFpfUmEQ'
$ws.Cells.Item(348, 2).Value = 'How to Repair the Durable Water Repellent (DWR) Finish on Your Rain Jacket'
$ws.Cells.Item(348, 3).Value = 'DIY'
$ws.Cells.Item(348, 4).Value = 'Other'
$ws.Cells.Item(348, 5).Value = '13:55'

# Row 349
$ws.Cells.Item(349, 1).Value = 'PF_WTo-uLvg'
$ws.Cells.Item(349, 2).Value = 'DIY Knife Sharpness Tester'
$ws.Cells.Item(349, 3).Value = 'DIY'
$ws.Cells.Item(349, 4).Value = 'Other'
$ws.Cells.Item(349, 5).Value = '8:45'

# Row 350
$ws.Cells.Item(350, 1).Value = 'fJMZkE6UxiI'
$ws.Cells.Item(350, 2).Value = 'How to Remove the Brake Cable on an Electric Scooter'
$ws.Cells.Item(350, 3).Value = 'DIY'
$ws.Cells.Item(350, 4).Value = 'Other'
$ws.Cells.Item(350, 5).Value = '3:17'

# Row 351
$ws.Cells.Item(351, 1).Value = 'tYkIt16bggw'
$ws.Cells.Item(351, 2).Value = 'Importing Table Data Into Matlab Using â€˜readtableâ€™'
$ws.Cells.Item(351, 3).Value = 'Matlab'
$ws.Cells.Item(351, 4).Value = 'Education'
$ws.Cells.Item(351, 5).Value = '30:19'

# Row 352
$ws.Cells.Item(352, 1).Value = '9gRTenEbwSk'
$ws.Cells.Item(352, 2).Value = 'Chefâ€™s Choice Model 15XV Electric Knife Sharpener: a Review with Quantitative Data'
$ws.Cells.Item(352, 3).Value = 'Review'
$ws.Cells.Item(352, 4).Value = 'Other'
$ws.Cells.Item(352, 5).Value = '19:14'

# Row 353
$ws.Cells.Item(353, 1).Value = 'XZt1YLVy6XU'
$ws.Cells.Item(353, 2).Value = 'AA516/AE512 Week02'
$ws.Cells.Item(353, 3).Value = 'FlightMechanics'
$ws.Cells.Item(353, 4).Value = 'Education'
$ws.Cells.Item(353, 5).Value = '10:15'

# Row 354
$ws.Cells.Item(354, 1).Value = 'NI70-AWnO4w'
$ws.Cells.Item(354, 2).Value = 'Direction Cosine Matrix from North East Down to East North Up'
$ws.Cells.Item(354, 3).Value = 'FlightMechanics'
$ws.Cells.Item(354, 4).Value = 'Education'
$ws.Cells.Item(354, 5).Value = '16:33'

# Row 355
$ws.Cells.Item(355, 1).Value = 'iaM-jMmAp6c'
$ws.Cells.Item(355, 2).Value = 'AA516/AE512 Week03'
$ws.Cells.Item(355, 3).Value = 'FlightMechanics'
$ws.Cells.Item(355, 4).Value = 'Education'
$ws.Cells.Item(355, 5).Value = '9:01'

# Row 356
$ws.Cells.Item(356, 1).Value = 'WDMTaNsgYVo'
$ws.Cells.Item(356, 2).Value = 'AA516/AE512 Week04'
$ws.Cells.Item(356, 3).Value = 'FlightMechanics'
$ws.Cells.Item(356, 4).Value = 'Education'
$ws.Cells.Item(356, 5).Value = '9:10'

# Row 357
$ws.Cells.Item(357, 1).Value = 'AlTYdT7kF38'
$ws.Cells.Item(357, 2).Value = 'Find the Radius of the Circle Inside the Triangle'
$ws.Cells.Item(357, 3).Value = 'FlightMechanics'
$ws.Cells.Item(357, 4).Value = 'Education'
$ws.Cells.Item(357, 5).Value = '13:13'

# Row 358
$ws.Cells.Item(358, 1).Value = 'js0jPF_h0wM'
$ws.Cells.Item(358, 2).Value = 'AA516/AE512 Week05'
$ws.Cells.Item(358, 3).Value = 'FlightMechanics'
$ws.Cells.Item(358, 4).Value = 'Education'
$ws.Cells.Item(358, 5).Value = '12:55'

# Row 359
$ws.Cells.Item(359, 1).Value = ' -FQzeD9gsS0'
$ws.Cells.Item(359, 2).Value = 'AA516/AE512 Week06'
$ws.Cells.Item(359, 3).Value = 'FlightMechanics'
$ws.Cells.Item(359, 4).Value = 'Education'
$ws.Cells.Item(359, 5).Value = '16:26'

# Row 360
$ws.Cells.Item(360, 1).Value = 'bJlryGNGDrw'
$ws.Cells.Item(360, 2).Value = 'MathWorks Free Stuff'
$ws.Cells.Item(360, 3).Value = 'Matlab'
$ws.Cells.Item(360, 4).Value = 'Other'
$ws.Cells.Item(360, 5).Value = '6:44'

# Row 361
$ws.Cells.Item(361, 1).Value = 'aDWz_hbNYzM'
$ws.Cells.Item(361, 2).Value = 'AA516/AE512 Week07'
$ws.Cells.Item(361, 3).Value = 'FlightMechanics'
$ws.Cells.Item(361, 4).Value = 'Education'
$ws.Cells.Item(361, 5).Value = '15:02'

# Row 362
$ws.Cells.Item(362, 1).Value = 'EGmKPpV6bXs'
$ws.Cells.Item(362, 2).Value = 'AA516/AE512 Week08'
$ws.Cells.Item(362, 3).Value = 'FlightMechanics'
$ws.Cells.Item(362, 4).Value = 'Education'
$ws.Cells.Item(362, 5).Value = '10:53'

# Row 363
$ws.Cells.Item(363, 1).Value = 'JZNORkEnB5c'
$ws.Cells.Item(363, 2).Value = 'AA516/AE512 Week09'
$ws.Cells.Item(363, 3).Value = 'FlightMechanics'
$ws.Cells.Item(363, 4).Value = 'Education'
$ws.Cells.Item(363, 5).Value = '19:14'

# Row 364
$ws.Cells.Item(364, 1).Value = 'Z5sWyYpLbnQ'
$ws.Cells.Item(364, 2).Value = 'AA516/AE512 Week10'
$ws.Cells.Item(364, 3).Value = 'FlightMechanics'
$ws.Cells.Item(364, 4).Value = 'Education'
$ws.Cells.Item(364, 5).Value = '9:48'

# Row 365
$ws.Cells.Item(365, 1).Value = 'n2bvoNsyhcE'
$ws.Cells.Item(365, 2).Value = 'REI Down Time 25 Sleeping Bag'
$ws.Cells.Item(365, 3).Value = 'Review'
$ws.Cells.Item(365, 4).Value = 'Other'
$ws.Cells.Item(365, 5).Value = '9:29'

# Row 366
$ws.Cells.Item(366, 1).Value = 'GLTenrOMBz8'
$ws.Cells.Item(366, 2).Value = 'Replace the Line on a DeWalt String Trimmer'
$ws.Cells.Item(366, 3).Value = 'DIY'
$ws.Cells.Item(366, 4).Value = 'Other'
$ws.Cells.Item(366, 5).Value = '10:03'

# Row 367
$ws.Cells.Item(367, 1).Value = 'w4mIw3kSTVU'
$ws.Cells.Item(367, 2).Value = 'Fuel Leaking from Motorcycle â€“ Solved!'
$ws.Cells.Item(367, 3).Value = 'DIY'
$ws.Cells.Item(367, 4).Value = 'Other'
$ws.Cells.Item(367, 5).Value = '5:45'

# Row 368
$ws.Cells.Item(368, 1).Value = 'x6qkmxxMKEw'
$ws.Cells.Item(368, 2).Value = 'Commenting Code in Matlab and Simulink (Ways You Can Use â€˜%â€™)'
$ws.Cells.Item(368, 3).Value = 'Matlab'
$ws.Cells.Item(368, 4).Value = 'Education'
$ws.Cells.Item(368, 5).Value = '29:50'

# Row 369
$ws.Cells.Item(369, 1).Value = 'm46opXldvEA'
$ws.Cells.Item(369, 2).Value = 'Disassemble the Headset and Top Tube on an Electric Scooter (GoTrax Eclipse)'
$ws.Cells.Item(369, 3).Value = 'DIY'
$ws.Cells.Item(369, 4).Value = 'Other'
$ws.Cells.Item(369, 5).Value = '5:31'

# Row 370
$ws.Cells.Item(370, 1).Value = 'NF71HPAB2W0'
$ws.Cells.Item(370, 2).Value = 'Default Arguments in Matlab Functions (varargin and nargin)'
$ws.Cells.Item(370, 3).Value = 'Matlab'
$ws.Cells.Item(370, 4).Value = 'Education'
$ws.Cells.Item(370, 5).Value = '17:34'

# Reflect the author's final cursor position after entering the new rows
$null = $ws.Range("F326").Select()

